# Refresh cached Universalis market-price snapshot values and derived profit
# figures on each job sheet (ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR) of the Excalibur
# Profits workbook, as produced by the scheduled data-pull runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 3494.3333
$ws.Range("I32").Value = 3799.5
$ws.Range("J32").Value = 3407.1428
$ws.Range("K32").Value = 3799.5
$ws.Range("L32").Value = 3407.1428
$ws.Range("M32").Value = -3473.5
$ws.Range("N32").Value = -4059.1428
$ws.Range("H34").Value = 5368.5
$ws.Range("I34").Value = 5368.5
$ws.Range("K34").Value = 5368.5
$ws.Range("M34").Value = -5165.5
$ws.Range("H36").Value = 5368.5
$ws.Range("I36").Value = 5368.5
$ws.Range("K36").Value = 5368.5
$ws.Range("M36").Value = -4653.5
$ws.Range("H40").Value = 1928.3334
$ws.Range("I40").Value = 1879.2222
$ws.Range("K40").Value = 1879.2222
$ws.Range("M40").Value = -1704.2222
$ws.Range("H64").Value = 6000
$ws.Range("H67").Value = 6000
$ws.Range("H111").Value = 3538.4707
$ws.Range("I111").Value = 1536
$ws.Range("J111").Value = 4630.727
$ws.Range("K111").Value = 4608
$ws.Range("L111").Value = 13892.181
$ws.Range("M111").Value = -1541
$ws.Range("N111").Value = -20026.181
$ws.Range("H112").Value = 4909.143
$ws.Range("J112").Value = 5072.239
$ws.Range("L112").Value = 15216.717
$ws.Range("N112").Value = -17432.717
$ws.Range("H132").Value = 139236
$ws.Range("I132").Value = 168830.58
$ws.Range("K132").Value = 506491.74
$ws.Range("M132").Value = -503961.74
$ws.Range("H141").Value = 1652.579
$ws.Range("I141").Value = 1688.8889
$ws.Range("K141").Value = 5066.6667
$ws.Range("M141").Value = 113.3333000000002

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H37").Value = 68505.625
$ws.Range("J37").Value = 77007.62
$ws.Range("L37").Value = 77007.62
$ws.Range("N37").Value = -77553.62
$ws.Range("H57").Value = 0
$ws.Range("I57").Value = 0
$ws.Range("K57").Value = 0
$ws.Range("H61").Value = 3705683.5
$ws.Range("I61").Value = 4763164.5
$ws.Range("K61").Value = 4763164.5
$ws.Range("M61").Value = -4762952.5
$ws.Range("H80").Value = 89388.8
$ws.Range("I80").Value = 75000
$ws.Range("J80").Value = 92986
$ws.Range("K80").Value = 75000
$ws.Range("L80").Value = 92986
$ws.Range("M80").Value = -74002
$ws.Range("N80").Value = -94982
$ws.Range("H83").Value = 89388.8
$ws.Range("I83").Value = 75000
$ws.Range("J83").Value = 92986
$ws.Range("K83").Value = 225000
$ws.Range("L83").Value = 278958
$ws.Range("M83").Value = -220008
$ws.Range("N83").Value = -288942
$ws.Range("H110").Value = 2046.75
$ws.Range("I110").Value = 1208
$ws.Range("K110").Value = 1208
$ws.Range("M110").Value = 837
$ws.Range("H136").Value = 3705683.5
$ws.Range("I136").Value = 4763164.5
$ws.Range("K136").Value = 14289493.5
$ws.Range("M136").Value = -14286943.5
$ws.Range("M57").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 2356
$ws.Range("I64").Value = 340
$ws.Range("J64").Value = 3700
$ws.Range("K64").Value = 340
$ws.Range("L64").Value = 3700
$ws.Range("M64").Value = -115
$ws.Range("N64").Value = -4150
$ws.Range("H67").Value = 2356
$ws.Range("I67").Value = 340
$ws.Range("J67").Value = 3700
$ws.Range("K67").Value = 340
$ws.Range("L67").Value = 3700
$ws.Range("M67").Value = 440
$ws.Range("N67").Value = -5260
$ws.Range("H86").Value = 1133.1666
$ws.Range("I86").Value = 1159.8
$ws.Range("K86").Value = 1159.8
$ws.Range("M86").Value = -36.79999999999995
$ws.Range("H89").Value = 1133.1666
$ws.Range("I89").Value = 1159.8
$ws.Range("K89").Value = 5799
$ws.Range("M89").Value = -183
$ws.Range("H99").Value = 27795.084
$ws.Range("I99").Value = 44901.918
$ws.Range("J99").Value = 10688.25
$ws.Range("K99").Value = 44901.918
$ws.Range("L99").Value = 10688.25
$ws.Range("M99").Value = -43403.918
$ws.Range("N99").Value = -13684.25
$ws.Range("H105").Value = 1633
$ws.Range("I105").Value = 1654.9231
$ws.Range("J105").Value = 1348
$ws.Range("K105").Value = 1654.9231
$ws.Range("L105").Value = 1348
$ws.Range("M105").Value = 92.07690000000002
$ws.Range("N105").Value = -4842

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 26223.105
$ws.Range("I31").Value = 11864.889
$ws.Range("J31").Value = 39145.5
$ws.Range("K31").Value = 11864.889
$ws.Range("L31").Value = 39145.5
$ws.Range("M31").Value = -11569.889
$ws.Range("N31").Value = -39735.5
$ws.Range("H34").Value = 26223.105
$ws.Range("I34").Value = 11864.889
$ws.Range("J34").Value = 39145.5
$ws.Range("K34").Value = 11864.889
$ws.Range("L34").Value = 39145.5
$ws.Range("M34").Value = -11662.889
$ws.Range("N34").Value = -39549.5
$ws.Range("H58").Value = 731466.4
$ws.Range("I58").Value = 1768138.2
$ws.Range("J58").Value = 5796
$ws.Range("K58").Value = 1768138.2
$ws.Range("L58").Value = 5796
$ws.Range("M58").Value = -1767935.2
$ws.Range("N58").Value = -6202
$ws.Range("H74").Value = 59971
$ws.Range("J74").Value = 59971
$ws.Range("L74").Value = 59971
$ws.Range("N74").Value = -61719
$ws.Range("H77").Value = 59971
$ws.Range("J77").Value = 59971
$ws.Range("L77").Value = 179913
$ws.Range("N77").Value = -188649
$ws.Range("H122").Value = 6730
$ws.Range("I122").Value = 3233.875
$ws.Range("J122").Value = 10226.125
$ws.Range("K122").Value = 9701.625
$ws.Range("L122").Value = 30678.375
$ws.Range("M122").Value = -7251.625
$ws.Range("N122").Value = -35578.375
$ws.Range("H134").Value = 9247.421
$ws.Range("I134").Value = 10796.733
$ws.Range("K134").Value = 32390.199
$ws.Range("M134").Value = -29855.199
$ws.Range("H136").Value = 731466.4
$ws.Range("I136").Value = 1768138.2
$ws.Range("J136").Value = 5796
$ws.Range("K136").Value = 5304414.6
$ws.Range("L136").Value = 17388
$ws.Range("M136").Value = -5301864.6
$ws.Range("N136").Value = -22488

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 2724.5
$ws.Range("J22").Value = 2724.5
$ws.Range("L22").Value = 8173.5
$ws.Range("N22").Value = -8511.5
$ws.Range("H27").Value = 2724.5
$ws.Range("J27").Value = 2724.5
$ws.Range("L27").Value = 8173.5
$ws.Range("N27").Value = -8377.5
$ws.Range("H75").Value = 5914.5
$ws.Range("I75").Value = 1770
$ws.Range("K75").Value = 5310
$ws.Range("M75").Value = -4312
$ws.Range("H78").Value = 5914.5
$ws.Range("I78").Value = 1770
$ws.Range("K78").Value = 15930
$ws.Range("M78").Value = -10938
$ws.Range("H110").Value = 1863.5
$ws.Range("I110").Value = 1863.5
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 5590.5
$ws.Range("L110").Value = 0
$ws.Range("M110").Value = -1500.5
$ws.Range("N110").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1564.0952
$ws.Range("I97").Value = 1574.8334
$ws.Range("J97").Value = 1499.6666
$ws.Range("K97").Value = 1574.8334
$ws.Range("L97").Value = 1499.6666
$ws.Range("M97").Value = -1078.8334
$ws.Range("N97").Value = -2491.6666
$ws.Range("H102").Value = 2874.2727
$ws.Range("J102").Value = 4509.875
$ws.Range("L102").Value = 4509.875
$ws.Range("N102").Value = -7753.875
$ws.Range("H113").Value = 5335.484
$ws.Range("I113").Value = 4458.3335
$ws.Range("K113").Value = 4458.3335
$ws.Range("M113").Value = -2288.3335
$ws.Range("H122").Value = 4421.8125
$ws.Range("I122").Value = 2203.9565
$ws.Range("K122").Value = 6611.869499999999
$ws.Range("M122").Value = -4161.869499999999
$ws.Range("H126").Value = 697771.4399999999
$ws.Range("I126").Value = 1193051.1
$ws.Range("K126").Value = 3579153.3
$ws.Range("M126").Value = -3576683.3
$ws.Range("H132").Value = 50608200
$ws.Range("I132").Value = 72293190
$ws.Range("K132").Value = 216879570
$ws.Range("M132").Value = -216877040

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4019.9
$ws.Range("I7").Value = 3733.2666
$ws.Range("K7").Value = 3733.2666
$ws.Range("M7").Value = -3621.2666
$ws.Range("H40").Value = 4150.5
$ws.Range("I40").Value = 4150.5
$ws.Range("K40").Value = 4150.5
$ws.Range("M40").Value = -4014.5
$ws.Range("H55").Value = 206.44444
$ws.Range("I55").Value = 192.75
$ws.Range("K55").Value = 192.75
$ws.Range("M55").Value = -19.75
$ws.Range("H61").Value = 1337.1666
$ws.Range("I61").Value = 1337.1666
$ws.Range("K61").Value = 1337.1666
$ws.Range("M61").Value = -1135.1666
$ws.Range("H113").Value = 1337.1666
$ws.Range("I113").Value = 1337.1666
$ws.Range("K113").Value = 1337.1666
$ws.Range("M113").Value = 832.8334
$ws.Range("H126").Value = 4019.9
$ws.Range("I126").Value = 3733.2666
$ws.Range("K126").Value = 11199.7998
$ws.Range("M126").Value = -8729.799800000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 805.5714
$ws.Range("I100").Value = 582.1875
$ws.Range("J100").Value = 1520.4
$ws.Range("K100").Value = 1164.375
$ws.Range("L100").Value = 3040.8
$ws.Range("M100").Value = -623.375
$ws.Range("N100").Value = -4122.8
$ws.Range("H113").Value = 2780.4707
$ws.Range("I113").Value = 550.5
$ws.Range("K113").Value = 1651.5
$ws.Range("M113").Value = 518.5

